# Update crypto price/volume table per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking price strings to stay as text (matches source data which
# stores prices like "226.08" or "0.547" as literal text, not numeric values).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated values
$ws.Range('D2').Value = '34.129.08'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '1.784.46'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').Value = '226.08'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').Value = '0.547'
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').Value = '31.86'
$ws.Range('E8').Value = '  -3.49%  '
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').Value = '0.0688'
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('D11').Value = '0.0943'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').Value = '2.041.60'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '11.19'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('D14').Value = '1.788.25'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('D15').Value = '34.025.90'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '67.92'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').Value = '245.72'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '0.0₃0778'
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').Value = '10.79'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('E24').Value = '  -1.73%  '
$ws.Range('D25').Value = '161.43'
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('D26').Value = '7.12'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').Value = '16.29'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = '3.65'
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('D33').Value = '3.59'
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('D34').Value = '1.81'
$ws.Range('E34').Value = '  -0.74%  '
$ws.Range('D35').Value = '1.461.77'
$ws.Range('E35').Value = '  +4.50%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.643'
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.0193'
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('D38').Value = '2.39'
$ws.Range('E38').Value = '  +7.75%  '
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').Value = '79.76'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').Value = '0.916'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').Value = '2.67'
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('D44').Value = '13.31'
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').Value = '0.0508'
$ws.Range('E45').Value = '  +2.47%  '
$ws.Range('D46').Value = '6.02'
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = '1.07'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0136'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('D49').Value = '107.15'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('D50').Value = '1.943.63'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('E51').Value = '  +0.62%  '
